$d = $word.ActiveDocument

# The contact-info paragraph (Block Text style) is the last paragraph in the
# document body; remember its index since we will append a new paragraph
# right after it.
$paraIndex = $d.Paragraphs.Count

# --- Remove everything that used to follow the hyperlink ------------------
# (" ", "• 215-983-0808 • 19 years old", the line break and the old address
# line). Doing this before inserting the new lead-in text keeps the new
# runs from being coalesced with anything else by the later edits.
$hl = $d.Hyperlinks.Item(1)
$hlRange = $hl.Range
$p = $d.Paragraphs.Item($paraIndex)
$tail = $d.Range($hlRange.End, $p.Range.End)
$tail.Text = ""

# --- Insert the new lead-in text right before the mailto hyperlink --------
# Insert the single-space run first (touching the hyperlink), then insert
# the longer "<<++>> ... 215-983-0808 ..." run right before *that* space so
# the two stay as distinct runs instead of being coalesced into one.
$hl = $d.Hyperlinks.Item(1)
$hlRange = $hl.Range
$insSpace = $d.Range($hlRange.Start, $hlRange.Start)
$insSpace.InsertBefore(" ")

$hl = $d.Hyperlinks.Item(1)
$hlRange = $hl.Range
$insLead = $d.Range($hlRange.Start - 1, $hlRange.Start - 1)
$insLead.InsertBefore("<<++>> • 215-983-0808 •")

# --- Split the address onto its own paragraph ------------------------------
# (same BlockText style) with the updated wording.
$p = $d.Paragraphs.Item($paraIndex)
$p.Range.InsertParagraphAfter()
$addrPara = $d.Paragraphs.Item($paraIndex + 1)
$addrPara.Range.Text = "21 Elm Avenue, Rockledge, Pennsylvania. 19046"
